# Add the "United States" column and refresh the survey figures now that
# 85% of the data has been collected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("F1").Value = "United States"

# Updated figures for France / Germany / Spain / United Kingdom, plus the
# new (still empty) United States column.
$data = @(
    @(0.1,    0.172043010752688, 0.12, 0.335365853658537),
    @(0.2,    0.505376344086022, 0.48, 0.402439024390244),
    @(0.3,    0.408602150537634, 0.68, 0.75),
    @(0.1,    0.193548387096774, 0.32, 0.201219512195122),
    @(0.325,  0.462365591397849, 0.52, 0.414634146341463),
    @(0.05,   0.021505376344086, 0.06, 0.0914634146341463),
    @(0.075,  0.032258064516129, 0.02, 0.00609756097560976),
    @(0.1,    0.032258064516129, 0.1,  0.0731707317073171),
    @(0.325,  0.333333333333333, 0.3,  0.268292682926829),
    @(0.25,   0.139784946236559, 0.36, 0.298780487804878)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
